# This script updates column G ("K") on Sheet1 with newly-computed
# strikeout-based K values, replacing the previous Strike# derived values.
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
#  calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row => new K value, per the regenerated save_data.
$kValues = @{
    2  = 1
    3  = 9
    4  = 5
    5  = 7
    6  = 2
    7  = 5
    8  = 4
    9  = 10
    10 = 4
    11 = 1
    12 = 6
    13 = 6
    14 = 2
    15 = 7
    16 = 3
    17 = 5
    18 = 8
    19 = 3
    20 = 3
    21 = 5
    22 = 3
    23 = 4
    24 = 7
    25 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
